# In the Education section, the paragraph ending with the bold date
# "April 2021" is preceded by a run of whitespace-only runs:
#   [91 spaces][20 spaces][9 spaces][3 spaces][22 spaces] "April 2021"
# The edit removes the two short runs (9 spaces + 3 spaces = 12
# characters total) that sit right before the final 22-space run, so
# the run sequence becomes:
#   [91 spaces][20 spaces][22 spaces] "April 2021"
# (the 22-space run itself is left untouched).

$d = $word.ActiveDocument

# "April 2021" is unique in the document, so use it as an anchor and
# compute the target range from fixed, content-intrinsic offsets: the
# surviving run is the 22 characters immediately before "April 2021",
# and the 12 characters to remove sit immediately before that run.
$anchor = $d.Content
$anchor.Find.Execute("April 2021", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$aprilStart = $anchor.Start
$keepRunStart = $aprilStart - 22
$removeStart = $keepRunStart - 12
$removeEnd = $keepRunStart

$toRemove = $d.Range($removeStart, $removeEnd)

# Use Find/Replace (rather than Range.Delete / Range.Text = "") so
# the untouched neighboring runs keep their own separate identities
# instead of being coalesced into a single run.
$toRemove.Find.Execute("            ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
